$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.211.99"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "3.941.29"
$ws.Range("E3").Value = "  +4.16%  "
$ws.Range("E4").Value = "  -0.06%  "
$s_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "471.59"
$ws.Range("D5").Style = $s_D5
$ws.Range("E5").Value = "  +9.31%  "
$s_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.17"
$ws.Range("D6").Style = $s_D6
$ws.Range("E6").Value = "  +4.27%  "
$s_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").Style = $s_D7
$ws.Range("E7").Value = "  +1.10%  "
$s_D8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("D8").Style = $s_D8
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +7.66%  "
$ws.Range("E11").Value = "  +6.85%  "
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "4.566.93"
$ws.Range("E13").Value = "  +4.21%  "
$s_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.36"
$ws.Range("D14").Style = $s_D14
$ws.Range("E14").Value = "  -0.80%  "
$s_D15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.28"
$ws.Range("D15").Style = $s_D15
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("D16").Value = "3.935.76"
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").Value = "67.512.00"
$ws.Range("E20").Value = "  +1.23%  "
$s_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "437.83"
$ws.Range("D21").Style = $s_D21
$ws.Range("E21").Value = "  +6.94%  "
$ws.Range("E22").Value = "  +4.87%  "
$s_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.52"
$ws.Range("D23").Style = $s_D23
$ws.Range("E23").Value = "  -1.38%  "
$s_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.73"
$ws.Range("D24").Style = $s_D24
$ws.Range("E24").Value = "  +2.81%  "
$s_D25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.63"
$ws.Range("D25").Style = $s_D25
$ws.Range("E25").Value = "  +8.39%  "
$ws.Range("E26").Value = "  +6.22%  "
$s_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.77"
$ws.Range("D27").Style = $s_D27
$ws.Range("E27").Value = "  +2.20%  "
$s_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.21"
$ws.Range("D28").Style = $s_D28
$ws.Range("E28").Value = "  +4.70%  "
$s_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.79"
$ws.Range("D29").Style = $s_D29
$ws.Range("E29").Value = "  +1.70%  "
$s_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "722.08"
$ws.Range("D30").Style = $s_D30
$ws.Range("E30").Value = "  +0.81%  "
$s_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.57"
$ws.Range("D31").Style = $s_D31
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("E32").Value = "  -2.11%  "
$s_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.78"
$ws.Range("D33").Style = $s_D33
$ws.Range("E33").Value = "  +3.32%  "
$s_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.71"
$ws.Range("D34").Style = $s_D34
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("E35").Value = "  +1.52%  "
$s_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.87"
$ws.Range("D36").Style = $s_D36
$ws.Range("E36").Value = "  +3.22%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "0.0₃0787"
$ws.Range("E38").Value = "  +16.56%  "
$s_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.37"
$ws.Range("D39").Style = $s_D39
$ws.Range("E39").Value = "  -5.89%  "
$ws.Range("E40").Value = "  +0.71%  "
$s_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.05"
$ws.Range("D41").Style = $s_D41
$ws.Range("E41").Value = "  +4.26%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$s_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.142"
$ws.Range("D42").Style = $s_D42
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$s_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.336"
$ws.Range("D43").Style = $s_D43
$ws.Range("E43").Value = "  +4.59%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$s_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.57"
$ws.Range("D44").Style = $s_D44
$ws.Range("E44").Value = "  -7.53%  "
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$s_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.49"
$ws.Range("D46").Style = $s_D46
$ws.Range("E46").Value = "  +4.73%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$s_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.80"
$ws.Range("D47").Style = $s_D47
$ws.Range("E47").Value = "  +4.03%  "
$ws.Range("E48").Value = "  +5.00%  "
$s_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "146.90"
$ws.Range("D49").Style = $s_D49
$ws.Range("E49").Value = "  +2.99%  "
$s_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.14"
$ws.Range("D50").Style = $s_D50
$ws.Range("E50").Value = "  -3.70%  "
$s_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.86"
$ws.Range("D51").Style = $s_D51
$ws.Range("E51").Value = "  +1.60%  "
